$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 23.061728
$ws.Cells.Item(2, 8).Value = 69.18518399999999
$ws.Cells.Item(2, 9).Value = 0.6130144106248721
$ws.Cells.Item(2, 10).Value = 0.6902769593117909
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 13.10301133333333
$ws.Cells.Item(2, 14).Value = 39.309034
$ws.Cells.Item(2, 15).Value = 0.2036753411238414
$ws.Cells.Item(2, 16).Value = 0.2126567538160978
$ws.Cells.Item(2, 17).Value = 302.1780833502507
$ws.Cells.Item(2, 18).Value = 2719.602750152256
$ws.Cells.Item(2, 19).Value = 0.1248559191978514
$ws.Cells.Item(2, 20).Value = 0.1467920574012921

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 23.061728
$ws.Cells.Item(3, 8).Value = 69.18518399999999
$ws.Cells.Item(3, 9).Value = 0.6130144106248721
$ws.Cells.Item(3, 10).Value = 0.6902769593117909
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 29.54200233333333
$ws.Cells.Item(3, 14).Value = 88.626007
$ws.Cells.Item(3, 15).Value = 0.4592056932299318
$ws.Cells.Item(3, 16).Value = 0.479455154057023
$ws.Cells.Item(3, 17).Value = 681.2896223866986
$ws.Cells.Item(3, 18).Value = 6131.606601480287
$ws.Cells.Item(3, 19).Value = 0.2814997073909324
$ws.Cells.Item(3, 20).Value = 0.3309568458688482

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 23.061728
$ws.Cells.Item(4, 8).Value = 69.18518399999999
$ws.Cells.Item(4, 9).Value = 0.6130144106248721
$ws.Cells.Item(4, 10).Value = 0.6902769593117909
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 5.915846333333334
$ws.Cells.Item(4, 14).Value = 17.747539
$ws.Cells.Item(4, 15).Value = 0.0919568784095198
$ws.Cells.Item(4, 16).Value = 0.09601187431786275
$ws.Cells.Item(4, 17).Value = 136.4296390291307
$ws.Cells.Item(4, 18).Value = 1227.866751262176
$ws.Cells.Item(4, 19).Value = 0.0563708916211148
$ws.Cells.Item(4, 20).Value = 0.06627478466196013

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 23.061728
$ws.Cells.Item(5, 8).Value = 69.18518399999999
$ws.Cells.Item(5, 9).Value = 0.6130144106248721
$ws.Cells.Item(5, 10).Value = 0.6902769593117909
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 7.620811666666666
$ws.Cells.Item(5, 14).Value = 22.862435
$ws.Cells.Item(5, 15).Value = 0.1184591370916581
$ws.Cells.Item(5, 16).Value = 0.123682795446755
$ws.Cells.Item(5, 17).Value = 175.7490857958933
$ws.Cells.Item(5, 18).Value = 1581.74177216304
$ws.Cells.Item(5, 19).Value = 0.07261715810737374
$ws.Cells.Item(5, 20).Value = 0.08537538396016824

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 23.061728
$ws.Cells.Item(6, 8).Value = 69.18518399999999
$ws.Cells.Item(6, 9).Value = 0.6130144106248721
$ws.Cells.Item(6, 10).Value = 0.6902769593117909
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 8.1511595
$ws.Cells.Item(6, 14).Value = 16.302319
$ws.Cells.Item(6, 15).Value = 0.1267029501450486
$ws.Cells.Item(6, 16).Value = 0.08819342236226138
$ws.Cells.Item(6, 17).Value = 187.979823273616
$ws.Cells.Item(6, 18).Value = 1127.878939641696
$ws.Cells.Item(6, 19).Value = 0.07767073430759955
$ws.Cells.Item(6, 20).Value = 0.06087788741952229

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1.926013333333334
$ws.Cells.Item(7, 8).Value = 5.778040000000001
$ws.Cells.Item(7, 9).Value = 0.05119624723650278
$ws.Cells.Item(7, 10).Value = 0.05764887294340218
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 13.10301133333333
$ws.Cells.Item(7, 14).Value = 39.309034
$ws.Cells.Item(7, 15).Value = 0.2036753411238414
$ws.Cells.Item(7, 16).Value = 0.2126567538160978
$ws.Cells.Item(7, 17).Value = 25.23657453481778
$ws.Cells.Item(7, 18).Value = 227.12917081336
$ws.Cells.Item(7, 19).Value = 0.01042741312015523
$ws.Cells.Item(7, 20).Value = 0.01225942218130058

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1.926013333333334
$ws.Cells.Item(8, 8).Value = 5.778040000000001
$ws.Cells.Item(8, 9).Value = 0.05119624723650278
$ws.Cells.Item(8, 10).Value = 0.05764887294340218
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 29.54200233333333
$ws.Cells.Item(8, 14).Value = 88.626007
$ws.Cells.Item(8, 15).Value = 0.4592056932299318
$ws.Cells.Item(8, 16).Value = 0.479455154057023
$ws.Cells.Item(8, 17).Value = 56.89829038736445
$ws.Cells.Item(8, 18).Value = 512.0846134862801
$ws.Cells.Item(8, 19).Value = 0.02350960820300924
$ws.Cells.Item(8, 20).Value = 0.02764004925829264

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1.926013333333334
$ws.Cells.Item(9, 8).Value = 5.778040000000001
$ws.Cells.Item(9, 9).Value = 0.05119624723650278
$ws.Cells.Item(9, 10).Value = 0.05764887294340218
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 5.915846333333334
$ws.Cells.Item(9, 14).Value = 17.747539
$ws.Cells.Item(9, 15).Value = 0.0919568784095198
$ws.Cells.Item(9, 16).Value = 0.09601187431786275
$ws.Cells.Item(9, 17).Value = 11.39399891595111
$ws.Cells.Item(9, 18).Value = 102.54599024356
$ws.Cells.Item(9, 19).Value = 0.004707847082150799
$ws.Cells.Item(9, 20).Value = 0.005534976343608368

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 1.926013333333334
$ws.Cells.Item(10, 8).Value = 5.778040000000001
$ws.Cells.Item(10, 9).Value = 0.05119624723650278
$ws.Cells.Item(10, 10).Value = 0.05764887294340218
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 7.620811666666666
$ws.Cells.Item(10, 14).Value = 22.862435
$ws.Cells.Item(10, 15).Value = 0.1184591370916581
$ws.Cells.Item(10, 16).Value = 0.123682795446755
$ws.Cells.Item(10, 17).Value = 14.67778488082222
$ws.Cells.Item(10, 18).Value = 132.1000639274
$ws.Cells.Item(10, 19).Value = 0.006064663269967307
$ws.Cells.Item(10, 20).Value = 0.007130173759994779

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 1.926013333333334
$ws.Cells.Item(11, 8).Value = 5.778040000000001
$ws.Cells.Item(11, 9).Value = 0.05119624723650278
$ws.Cells.Item(11, 10).Value = 0.05764887294340218
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 8.1511595
$ws.Cells.Item(11, 14).Value = 16.302319
$ws.Cells.Item(11, 15).Value = 0.1267029501450486
$ws.Cells.Item(11, 16).Value = 0.08819342236226138
$ws.Cells.Item(11, 17).Value = 15.69924187912667
$ws.Cells.Item(11, 18).Value = 94.19545127476002
$ws.Cells.Item(11, 19).Value = 0.006486715561220196
$ws.Cells.Item(11, 20).Value = 0.005084251400205811

$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 7).Value = 12.632464
$ws.Cells.Item(12, 8).Value = 25.264928
$ws.Cells.Item(12, 9).Value = 0.3357893421386252
$ws.Cells.Item(12, 10).Value = 0.2520741677448068
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 13.10301133333333
$ws.Cells.Item(12, 14).Value = 39.309034
$ws.Cells.Item(12, 15).Value = 0.2036753411238414
$ws.Cells.Item(12, 16).Value = 0.2126567538160978
$ws.Cells.Item(12, 17).Value = 165.5233189599253
$ws.Cells.Item(12, 18).Value = 993.139913759552
$ws.Cells.Item(12, 19).Value = 0.06839200880583479
$ws.Cells.Item(12, 20).Value = 0.05360527423350513

$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 7).Value = 12.632464
$ws.Cells.Item(13, 8).Value = 25.264928
$ws.Cells.Item(13, 9).Value = 0.3357893421386252
$ws.Cells.Item(13, 10).Value = 0.2520741677448068
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 29.54200233333333
$ws.Cells.Item(13, 14).Value = 88.626007
$ws.Cells.Item(13, 15).Value = 0.4592056932299318
$ws.Cells.Item(13, 16).Value = 0.479455154057023
$ws.Cells.Item(13, 17).Value = 373.1882809637493
$ws.Cells.Item(13, 18).Value = 2239.129685782496
$ws.Cells.Item(13, 19).Value = 0.1541963776359901
$ws.Cells.Item(13, 20).Value = 0.1208582589298822

$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 7).Value = 12.632464
$ws.Cells.Item(14, 8).Value = 25.264928
$ws.Cells.Item(14, 9).Value = 0.3357893421386252
$ws.Cells.Item(14, 10).Value = 0.2520741677448068
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 5.915846333333334
$ws.Cells.Item(14, 14).Value = 17.747539
$ws.Cells.Item(14, 15).Value = 0.0919568784095198
$ws.Cells.Item(14, 16).Value = 0.09601187431786275
$ws.Cells.Item(14, 17).Value = 74.73171583536534
$ws.Cells.Item(14, 18).Value = 448.3902950121919
$ws.Cells.Item(14, 19).Value = 0.0308781397062542
$ws.Cells.Item(14, 20).Value = 0.02420211331229425

$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 7).Value = 12.632464
$ws.Cells.Item(15, 8).Value = 25.264928
$ws.Cells.Item(15, 9).Value = 0.3357893421386252
$ws.Cells.Item(15, 10).Value = 0.2520741677448068
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 7.620811666666666
$ws.Cells.Item(15, 14).Value = 22.862435
$ws.Cells.Item(15, 15).Value = 0.1184591370916581
$ws.Cells.Item(15, 16).Value = 0.123682795446755
$ws.Cells.Item(15, 17).Value = 96.26962902994664
$ws.Cells.Item(15, 18).Value = 577.6177741796799
$ws.Cells.Item(15, 19).Value = 0.03977731571431711
$ws.Cells.Item(15, 20).Value = 0.03117723772659195

$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 7).Value = 12.632464
$ws.Cells.Item(16, 8).Value = 25.264928
$ws.Cells.Item(16, 9).Value = 0.3357893421386252
$ws.Cells.Item(16, 10).Value = 0.2520741677448068
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 8.1511595
$ws.Cells.Item(16, 14).Value = 16.302319
$ws.Cells.Item(16, 15).Value = 0.1267029501450486
$ws.Cells.Item(16, 16).Value = 0.08819342236226138
$ws.Cells.Item(16, 17).Value = 102.969228942008
$ws.Cells.Item(16, 18).Value = 411.876915768032
$ws.Cells.Item(16, 19).Value = 0.04254550027622891
$ws.Cells.Item(16, 20).Value = 0.02223128354253327

